$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.039.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.93%  '

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.17%  '

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.68%  '

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.39%  '

# Row 6: USDC -> USDC
$ws.Range("E6").Value = '  -0.73%  '

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4821'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.86%  '

# Row 8: Cardano -> Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3811'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.94%  '

# Row 9: Dogecoin -> Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07356'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '

# Row 10: Polygon -> Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9329'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.46%  '

# Row 11: Solana -> Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '

# Row 12: TRON -> TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07788'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.91%  '

# Row 13: WrappedEther -> WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.936.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.78%  '

# Row 14: Polkadot -> Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.501'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.16%  '

# Row 15: Chainlink -> Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.630'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.84%  '

# Row 16: Litecoin -> Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.09%  '

# Row 17: BinanceUSD -> BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.64%  '

# Row 18: ShibaInu -> ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008839'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.74%  '

# Row 19: Dai -> Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '

# Row 20: WrappedBTC -> WrappedBTC
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.073.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.04%  '

# Row 21: Avalanche -> Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.70%  '

# Row 22: Uniswap -> Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.180'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.96%  '

# Row 23: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.148.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.57%  '

# Row 24: Cosmos -> Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.78%  '

# Row 25: Monero -> Toncoin
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.922'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.75%  '

# Row 26: Toncoin -> Monero
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '

# Row 27: EthereumClassic -> EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.11%  '

# Row 28: LidoDAOToken -> LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.115'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.35%  '

# Row 29: BitcoinCash -> BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.55%  '

# Row 30: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.960'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.82%  '

# Row 31: Stellar -> Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08958'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.35%  '

# Row 32: HuobiToken -> HuobiToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.302'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.12%  '

# Row 33: ARBITRUM -> ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.251'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.94%  '

# Row 34: ImmutableX -> ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7752'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.92%  '

# Row 35: Filecoin -> Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.676'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.35%  '

# Row 36: RenderToken -> RenderToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.09%  '

# Row 37: VeChain -> VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02054'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.10%  '

# Row 38: TrustWalletToken -> TrustWalletToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.112'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '

# Row 39: Hedera -> Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05308'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.61%  '

# Row 40: TheSandbox -> TheSandbox
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5478'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.41%  '

# Row 41: MXToken -> MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.996'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.25%  '

# Row 42: FraxShare -> FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.023'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '

# Row 43: Aptos -> Algorand
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1527'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '

# Row 44: Algorand -> Aptos
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.503'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.25%  '

# Row 45: EnergySwap -> EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '

# Row 46: Decentraland -> Quant
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '108.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.37%  '

# Row 47: Quant -> Decentraland
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4832'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.50%  '

# Row 48: PaxDollar -> PaxDollar
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.76%  '

# Row 49: NEARProtocol -> NEARProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.654'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '

# Row 50: Aave -> Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.98%  '

